$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark from its old location (right after
#    "...personas sin ").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. "Modulo de seguimiento alimentario." -> "Modulo de diagnostico y
#    seguimiento alimentario."
# ------------------------------------------------------------------
$d.Content.Find.Execute("seguimiento alimentario.", $true, $false, $false, $false, $false, $true, 1, $false, "diagnóstico y seguimiento alimentario.", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Rewrite the paragraph describing the patient follow-up, splitting
#    it into: a paragraph about new patients, a blank paragraph, a
#    paragraph about returning patients and a paragraph introducing the
#    "El seguimiento alimentario incluye:" list.
# ------------------------------------------------------------------
$oldPara = "Se da un seguimiento de los pacientes a través de un historial médico, revaloraciones, historial alimentario y un diagnóstico clínico. El seguimiento alimentario incluye:"

$newPara1 = "En los pacientes nuevos se un primer diagnóstico que incluye: historia alimentaria, recordatorio de 24 horas, recordatorio de fin de semana, antecedentes heredo-familiares, encuesta de hábitos alimenticios y se realiza orientación alimentaria con el plato del buen comer."
$newPara2 = "En pacientes que no son nuevos se realiza un seguimiento alimentario a través de un historial médico de sus consultas; durante cada consulta se hacen revaloraciones, esto incluye un recordatorio de 24 horas, toma de medidas antropométricas y se asigna un plan de alimentación de acuerda a la necesidad o padecimiento."
$newPara3 = "El seguimiento alimentario incluye:"

$combined = $newPara1 + $newPara2 + $newPara3

$d.Content.Find.Execute($oldPara, $true, $false, $false, $false, $false, $true, 1, $false, $combined, 2) | Out-Null

# Split the combined text into three separate paragraphs, working from
# the end of the text backwards so earlier found ranges stay valid.
$r1 = $d.Content
$r1.Find.Execute($newPara3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Collapse(1)
$r1.InsertParagraphBefore()

$r2 = $d.Content
$r2.Find.Execute($newPara2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Collapse(1)
$r2.InsertParagraphBefore()

# Insert a genuinely empty paragraph between the "new patients" text and
# the "returning patients" text.
$r3 = $d.Content
$r3.Find.Execute($newPara2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Collapse(1)
$r3.InsertParagraphBefore()

# ------------------------------------------------------------------
# 4. "Estudio de un diagnostico clinico para detectar enfermedades que
#    guarden relacion." -> "...con relacion a la alimentacion."
# ------------------------------------------------------------------
$d.Content.Find.Execute("Estudio de un diagnóstico clínico para detectar enfermedades que guarden relación.", $true, $false, $false, $false, $false, $true, 1, $false, "Estudio de un diagnóstico clínico para detectar enfermedades con relación a la alimentación.", 2) | Out-Null

# ------------------------------------------------------------------
# 5. Re-touch the "Implementacion." run so any stale rendering cache
#    (lastRenderedPageBreak) tied to it is dropped.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Implementación.", $true, $false, $false, $false, $false, $true, 1, $false, "Implementación.", 2) | Out-Null

# ------------------------------------------------------------------
# 6. Add the "_GoBack" bookmark at its new location, right before the
#    "Pruebas." run.
# ------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Pruebas.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r4)

Write-Output "done"
